$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (values must remain as literal text, matching the
# original inlineStr cell type, so we force Text number format before assignment
# and then restore the default "Normal" style so no stray formatting is left behind.
$cellUpdates = @{
    "D2" = "287.50"
    "E2" = "1.18%"
    "G2" = "8"
    "D3" = "29.63"
    "E3" = "3.83%"
    "G3" = "8"
    "D4" = "5.091"
    "E4" = "0.34%"
    "G4" = "8"
    "D5" = "0.06692"
    "E5" = "3.53%"
    "G5" = "8"
    "D6" = "7.337"
    "E6" = "1.61%"
    "G6" = "8"
    "E7" = "1.17%"
    "G7" = "8"
    "D8" = "1.361"
    "E8" = "1.84%"
    "G8" = "8"
    "D9" = "0.9189"
    "E9" = "0.79%"
    "G9" = "8"
    "D10" = "0.1587"
    "E10" = "2.43%"
    "G10" = "8"
    "D11" = "0.06721"
    "E11" = "4.61%"
    "G11" = "8"
    "D12" = "0.07650"
    "E12" = "0.99%"
    "G12" = "8"
    "D13" = "0.02935"
    "E13" = "-1.78%"
    "G13" = "8"
    "D14" = "0.08980"
    "E14" = "0.29%"
    "G14" = "8"
    "D15" = "0.001577"
    "E15" = "-1.20%"
    "G15" = "8"
    "D16" = "0.04499"
    "E16" = "1.06%"
    "G16" = "8"
    "D17" = "0.0006478"
    "E17" = "-0.59%"
    "G17" = "8"
    "D18" = "0.006254"
    "E18" = "1.68%"
    "G18" = "8"
    "D19" = "3.438"
    "E19" = "-0.50%"
    "G19" = "8"
    "D20" = "2.227"
    "E20" = "-0.65%"
    "G20" = "8"
    "E21" = "0.75%"
    "G21" = "8"
    "E22" = "-2.43%"
    "G22" = "8"
    "D23" = "4.094"
    "E23" = "2.84%"
    "G23" = "8"
    "E24" = "1.65%"
    "G24" = "8"
    "D25" = "0.001188"
    "E25" = "0.30%"
    "G25" = "8"
    "E26" = "-4.60%"
    "G26" = "8"
    "E27" = "-0.18%"
    "G27" = "8"
    "D28" = "0.0001617"
    "E28" = "-1.28%"
    "G28" = "8"
    "G29" = "8"
    "G30" = "8"
    "G31" = "8"
    "G32" = "8"
    "G33" = "8"
    "G34" = "8"
    "G35" = "8"
    "G36" = "8"
    "G37" = "8"
    "G38" = "8"
    "G39" = "8"
    "D40" = "0.04272"
    "E40" = "3.30%"
    "G40" = "8"
    "D41" = "0.006722"
    "E41" = "-0.24%"
    "G41" = "8"
    "D42" = "0.1238"
    "E42" = "0.69%"
    "G42" = "8"
    "D43" = "0.002239"
    "E43" = "6.48%"
    "G43" = "8"
    "D44" = "0.01329"
    "E44" = "12.78%"
    "G44" = "8"
    "D45" = "0.00005701"
    "E45" = "5.71%"
    "G45" = "8"
    "D46" = "1.968"
    "E46" = "-3.59%"
    "G46" = "8"
    "D47" = "0.01306"
    "E47" = "-29.44%"
    "G47" = "8"
    "G48" = "8"
    "G49" = "8"
    "G50" = "8"
    "G51" = "8"
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
    $cell.Style = "Normal"
}
